$wb = $excel.ActiveWorkbook

# --- Flight Schedule sheet: insert a "Type" column before ETA/ETD ---
$ws1 = $wb.Worksheets.Item("Flight Schedule")
$ws1.Columns("D:D").Insert()
$ws1.Range("D1").Value = "Type"
$ws1.Range("D2").Value = "FSNC"
$ws1.Range("D3").Value = "FSNC"
$ws1.Range("D4").Value = "Low-Cost"
[void]$ws1.Range("D4").Select()

# --- Gates sheet: add a "Type" column, widen a compatible-AC entry, and add a new gate ---
$ws2 = $wb.Worksheets.Item("Gates")
$ws2.Range("D1").Value = "Type"
$ws2.Range("D2").Value = "Remote"
$ws2.Range("C3").Value = "B738, B3XM, E170"
$ws2.Range("D3").Value = "Jet Bridge"
$ws2.Range("D4").Value = "Jet Bridge"
$ws2.Range("A5").Value = "G4"
$ws2.Range("B5").Value = 400
$ws2.Range("C5").Value = "B738"
$ws2.Range("D5").Value = "Remote"
[void]$ws2.Range("C16").Select()
$ws2.Activate()
